# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gets a brand-new, still-empty column
# inserted right before the existing "Late" column (which, together with
# everything to its right, shifts one column over: N->O, O->P, P->Q).
# The new column inherits its width from the column immediately to its
# left (M). Finally, the "Repayment schedule" tab becomes the active
# sheet/tab (it was "Edit Repayment Schedule" before), with J14 selected
# on it.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Remember column M's width so the newly inserted column N can inherit it,
# the way Excel's own "Insert" does.
$leftColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# Insert a new blank column at N; existing N/O/P (Late / heading / Outstanding)
# shift right to O/P/Q.
$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = $leftColumnWidth

# Make "Repayment schedule" the active sheet/tab with J14 selected,
# which also clears the previous tab selection on "Edit Repayment Schedule".
$wsSchedule.Activate()
$wsSchedule.Range("J14").Select()
